$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new column D with model's saved feature values
$ws.Range("D1").Value = "hi"
$ws.Range("D5").Value = 1
$ws.Range("D7").Value = 1
